# Update the Monte-Carlo replication results on the "results" sheet
# with the refreshed simulation run (final production numbers).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("results")

$ws.Range("C2").Value = 0.059460625900946945
$ws.Range("E2").Value = 8.1790000000000003
$ws.Range("F2").Value = 176.68199999999999
$ws.Range("G2").Value = 58.933999999999997
$ws.Range("H2").Value = 8.3309999999999995
$ws.Range("I2").Value = 8.4469999999999992
$ws.Range("J2").Value = 0.058961693579078961
$ws.Range("L2").Value = 76.682000000000002
$ws.Range("M2").Value = 176.68199999999999
$ws.Range("N2").Value = 8.6389999999999993
$ws.Range("O2").Value = 8.4830000000000005

$ws.Range("C3").Value = 0.055361564810302695
$ws.Range("E3").Value = 8.1769999999999996
$ws.Range("F3").Value = 176.762
$ws.Range("G3").Value = 59.225999999999999
$ws.Range("H3").Value = 5.9930000000000003
$ws.Range("I3").Value = 5.907
$ws.Range("J3").Value = 0.060522602803304598
$ws.Range("L3").Value = 76.762
$ws.Range("M3").Value = 107.788
$ws.Range("N3").Value = 5.98
$ws.Range("O3").Value = 5.8579999999999997

$ws.Range("C4").Value = 0.059499458917413997
$ws.Range("E4").Value = 8.19
$ws.Range("F4").Value = 176.863
$ws.Range("G4").Value = 59.231000000000002
$ws.Range("H4").Value = 8.6069999999999993
$ws.Range("I4").Value = 8.19
$ws.Range("J4").Value = 0.059425665341082949
$ws.Range("L4").Value = 76.863
$ws.Range("M4").Value = 110.943
$ws.Range("N4").Value = 8.3960000000000008
$ws.Range("O4").Value = 8.23

$ws.Range("C5").Value = 0.05759178755211973
$ws.Range("E5").Value = 8.17
$ws.Range("F5").Value = 176.47300000000001
$ws.Range("G5").Value = 58.710999999999999
$ws.Range("H5").Value = 14.445
$ws.Range("I5").Value = 13.321
$ws.Range("J5").Value = 0.056571638674773217
$ws.Range("L5").Value = 76.472999999999999
$ws.Range("M5").Value = 117.964
$ws.Range("N5").Value = 13.858000000000001
$ws.Range("O5").Value = 13.593999999999999

$ws.Range("C6").Value = 0.063214954058388409
$ws.Range("E6").Value = 8.1920000000000002
$ws.Range("F6").Value = 178.27199999999999
$ws.Range("G6").Value = 61.677
$ws.Range("H6").Value = 8.6229999999999993
$ws.Range("I6").Value = 9.0730000000000004
$ws.Range("J6").Value = 0.058734806056735415
$ws.Range("L6").Value = 78.272000000000006
$ws.Range("M6").Value = 178.27199999999999
$ws.Range("N6").Value = 8.68
$ws.Range("O6").Value = 8.9819999999999993

$ws.Range("C7").Value = 0.063028857516570891
$ws.Range("E7").Value = 8.1829999999999998
$ws.Range("F7").Value = 178.54499999999999
$ws.Range("G7").Value = 61.725999999999999
$ws.Range("H7").Value = 6.1970000000000001
$ws.Range("I7").Value = 6.4130000000000003
$ws.Range("J7").Value = 0.058355795148247964
$ws.Range("L7").Value = 78.545000000000002
$ws.Range("M7").Value = 107.956
$ws.Range("N7").Value = 6.2679999999999998
$ws.Range("O7").Value = 6.16

$ws.Range("C8").Value = 0.061289491840871221
$ws.Range("E8").Value = 8.1869999999999994
$ws.Range("F8").Value = 178.51900000000001
$ws.Range("G8").Value = 61.796999999999997
$ws.Range("H8").Value = 8.9540000000000006
$ws.Range("I8").Value = 8.48
$ws.Range("J8").Value = 0.065825414135119747
$ws.Range("L8").Value = 78.519000000000005
$ws.Range("M8").Value = 111.37
$ws.Range("N8").Value = 9.0009999999999994
$ws.Range("O8").Value = 8.83

$ws.Range("C9").Value = 0.064804688505314012
$ws.Range("E9").Value = 8.1929999999999996
$ws.Range("F9").Value = 178.77799999999999
$ws.Range("G9").Value = 62.000999999999998
$ws.Range("H9").Value = 14.224
$ws.Range("I9").Value = 15.202999999999999
$ws.Range("J9").Value = 0.053070855334334574
$ws.Range("L9").Value = 78.778000000000006
$ws.Range("M9").Value = 118.583
$ws.Range("N9").Value = 14.619
$ws.Range("O9").Value = 14.374000000000001

$ws.Range("C10").Value = 0.078420665680586435
$ws.Range("E10").Value = 8.1690000000000005
$ws.Range("F10").Value = 181.86600000000001
$ws.Range("G10").Value = 67.067999999999998
$ws.Range("H10").Value = 9.7219999999999995
$ws.Range("I10").Value = 9.5570000000000004
$ws.Range("J10").Value = 0.077936820781322316
$ws.Range("L10").Value = 81.866
$ws.Range("M10").Value = 181.86600000000001
$ws.Range("N10").Value = 8.98
$ws.Range("O10").Value = 9.0820000000000007

$ws.Range("C11").Value = 0.079625837543733813
$ws.Range("E11").Value = 8.1950000000000003
$ws.Range("F11").Value = 181.80099999999999
$ws.Range("G11").Value = 67.028999999999996
$ws.Range("H11").Value = 6.7080000000000002
$ws.Range("I11").Value = 6.6219999999999999
$ws.Range("J11").Value = 0.087081624160275886
$ws.Range("L11").Value = 81.801000000000002
$ws.Range("M11").Value = 108.042
$ws.Range("N11").Value = 6.5750000000000002
$ws.Range("O11").Value = 6.4359999999999999

$ws.Range("C12").Value = 0.076547339998814384
$ws.Range("E12").Value = 8.1850000000000005
$ws.Range("F12").Value = 181.90799999999999
$ws.Range("G12").Value = 67.108000000000004
$ws.Range("H12").Value = 9.39
$ws.Range("I12").Value = 9.141
$ws.Range("J12").Value = 0.075076043633312384
$ws.Range("L12").Value = 81.908000000000001
$ws.Range("M12").Value = 111.813
$ws.Range("N12").Value = 9.66
$ws.Range("O12").Value = 9.4960000000000004

$ws.Range("C13").Value = 0.076419120428973211
$ws.Range("E13").Value = 8.1630000000000003
$ws.Range("F13").Value = 181.7
$ws.Range("G13").Value = 67.058999999999997
$ws.Range("H13").Value = 16.262
$ws.Range("I13").Value = 14.622
$ws.Range("J13").Value = 0.084649907340262759
$ws.Range("L13").Value = 81.7
$ws.Range("M13").Value = 119.11
$ws.Range("N13").Value = 15.617000000000001
$ws.Range("O13").Value = 15.279

$ws.Range("C14").Value = 0.099374222154186953
$ws.Range("E14").Value = 8.1880000000000006
$ws.Range("F14").Value = 183.453
$ws.Range("G14").Value = 69.73
$ws.Range("H14").Value = 9.94
$ws.Range("I14").Value = 9.69
$ws.Range("J14").Value = 0.098495500675524733
$ws.Range("L14").Value = 83.453000000000003
$ws.Range("M14").Value = 183.453
$ws.Range("N14").Value = 10.375
$ws.Range("O14").Value = 9.7129999999999992

$ws.Range("C15").Value = 0.10064718976926663
$ws.Range("E15").Value = 8.1829999999999998
$ws.Range("F15").Value = 183.739
$ws.Range("G15").Value = 69.876000000000005
$ws.Range("H15").Value = 7.0910000000000002
$ws.Range("I15").Value = 6.7759999999999998
$ws.Range("J15").Value = 0.09792746113989638
$ws.Range("L15").Value = 83.739000000000004
$ws.Range("M15").Value = 108.491
$ws.Range("N15").Value = 7.1669999999999998
$ws.Range("O15").Value = 7.0289999999999999

$ws.Range("C16").Value = 0.097702965483446619
$ws.Range("E16").Value = 8.17
$ws.Range("F16").Value = 183.42099999999999
$ws.Range("G16").Value = 69.841999999999999
$ws.Range("H16").Value = 10.664
$ws.Range("I16").Value = 10.23
$ws.Range("J16").Value = 0.094335572974990678
$ws.Range("L16").Value = 83.421000000000006
$ws.Range("M16").Value = 111.899
$ws.Range("N16").Value = 9.952
$ws.Range("O16").Value = 9.7530000000000001

$ws.Range("C17").Value = 0.097206417857266261
$ws.Range("E17").Value = 8.1649999999999991
$ws.Range("F17").Value = 183.55699999999999
$ws.Range("G17").Value = 69.894999999999996
$ws.Range("H17").Value = 16.459
$ws.Range("I17").Value = 16.268000000000001
$ws.Range("J17").Value = 0.10490814683361278
$ws.Range("L17").Value = 83.557000000000002
$ws.Range("M17").Value = 119.503
$ws.Range("N17").Value = 16.279
$ws.Range("O17").Value = 15.9
